# Correcion a Diebold Mariano y revision de Cap1
# Adds the "d=6" column and the corresponding ARMA_I(p,6,q) model rows that
# were missing from the MOD_3 EnCQR-LSTM results sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-LabelFormat($range) {
    $range.Font.Bold = $true
    $range.HorizontalAlignment = [Microsoft.Office.Interop.Excel.XlHAlign]::xlHAlignCenter
    $range.VerticalAlignment = [Microsoft.Office.Interop.Excel.XlVAlign]::xlVAlignTop
    $range.Borders.LineStyle = [Microsoft.Office.Interop.Excel.XlLineStyle]::xlContinuous
}

# ---------------------------------------------------------------------------
# 1) Insert a new column G for "d=6" (old column G "d=7" -> H, old H "d=10" -> I)
# ---------------------------------------------------------------------------
$ws.Columns("G:G").Insert()
$ws.Range("G1").Value = "d=6"

# ---------------------------------------------------------------------------
# 2) Insert the three new ARMA_I(0,6,q) rows, right after ARMA_I(0,5,2)
# ---------------------------------------------------------------------------
$ws.Rows("20:22").Insert()
Set-LabelFormat($ws.Range("A20:A22"))
$ws.Range("A20").Value = "ARMA_I(0,6,0)"
$ws.Range("G20").Value = 97.96121861706386
$ws.Range("A21").Value = "ARMA_I(0,6,1)"
$ws.Range("G21").Value = 97.89713001403216
$ws.Range("A22").Value = "ARMA_I(0,6,2)"
$ws.Range("G22").Value = 98.05770385741084

# ---------------------------------------------------------------------------
# 3) Insert the two new ARMA_I(1,6,q) rows, right after ARMA_I(1,5,1)
# ---------------------------------------------------------------------------
$ws.Rows("38:39").Insert()
Set-LabelFormat($ws.Range("A38:A39"))
$ws.Range("A38").Value = "ARMA_I(1,6,0)"
$ws.Range("G38").Value = 98.05577848879807
$ws.Range("A39").Value = "ARMA_I(1,6,1)"
$ws.Range("G39").Value = 97.86090428190013

# ---------------------------------------------------------------------------
# 4) Insert the two new ARMA_I(2,6,q) rows, right after ARMA_I(2,5,2)
# ---------------------------------------------------------------------------
$ws.Rows("54:55").Insert()
Set-LabelFormat($ws.Range("A54:A55"))
$ws.Range("A54").Value = "ARMA_I(2,6,0)"
$ws.Range("G54").Value = 98.00335357808103
$ws.Range("A55").Value = "ARMA_I(2,6,2)"
$ws.Range("G55").Value = 97.85146112799239
